# Revert the "special" time-label cells back to "time" (undoing the
# previous merge that swapped the shared-string value used by B7/F7/B12),
# and restore the prior selection (A3:B3) that was active before that edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B7 / F7 / B12 currently read "special" -- change them back to "time".
$ws.Range("B7").Value = "time"
$ws.Range("F7").Value = "time"
$ws.Range("B12").Value = "time"

# Restore the previously-active selection (A3:B3, anchored at A3).
$ws.Range("A3:B3").Select()
